$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SamplesTab query text (B3): the `Tumor` column now pulls
# straight from samp.sample_tumor_status instead of the collected `tumor`
# list gathered earlier in the query.
$newQuery = "MATCH (s:study)<--(p:participant)<--(samp:sample)`nWHERE s.study_name in [`"Washington University PDX Development and Trial Center`"]`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"

$ws.Range("B3").Value = $newQuery

# Move the active selection, matching the author's last interaction with
# the sheet before saving.
$ws.Range("B12").Select()
